# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 15032
$wsExhibition.Range("F3").Value = 19098
$wsExhibition.Range("F5").Value = 141
$wsExhibition.Range("F13").Value = 58
$wsExhibition.Range("F14").Value = 165
$wsExhibition.Range("F22").Value = 7971
$wsExhibition.Range("F29").Value = 6064
$wsExhibition.Range("F35").Value = 5449
$wsExhibition.Range("F36").Value = 565
$wsExhibition.Range("F37").Value = 15
$wsExhibition.Range("F38").Value = 29
$wsExhibition.Range("F39").Value = 48

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 15032
$wsAll.Range("F3").Value = 19098
$wsAll.Range("F5").Value = 141
$wsAll.Range("F13").Value = 58
$wsAll.Range("F14").Value = 165
$wsAll.Range("F23").Value = 7971
$wsAll.Range("F32").Value = 6064
$wsAll.Range("F38").Value = 5449
$wsAll.Range("F39").Value = 565
$wsAll.Range("F40").Value = 15
$wsAll.Range("F41").Value = 29
$wsAll.Range("F42").Value = 48
